$d = $word.ActiveDocument
$d.Content.Find.Execute("52+8=", $true, $true, $false, $false, $false, $true, 1, $false, "99-88=", 2) | Out-Null
$d.Content.Find.Execute("69-33=", $true, $true, $false, $false, $false, $true, 1, $false, "51-0=", 2) | Out-Null
$d.Content.Find.Execute("4+59=", $true, $true, $false, $false, $false, $true, 1, $false, "32+4=", 2) | Out-Null
$d.Content.Find.Execute("48+23=", $true, $true, $false, $false, $false, $true, 1, $false, "55-45=", 2) | Out-Null
$d.Content.Find.Execute("33+50=", $true, $true, $false, $false, $false, $true, 1, $false, "66-38=", 2) | Out-Null
$d.Content.Find.Execute("19+27=", $true, $true, $false, $false, $false, $true, 1, $false, "97-90=", 2) | Out-Null
$d.Content.Find.Execute("82-80=", $true, $true, $false, $false, $false, $true, 1, $false, "98-41=", 2) | Out-Null
$d.Content.Find.Execute("49+9=", $true, $true, $false, $false, $false, $true, 1, $false, "60-28=", 2) | Out-Null
$d.Content.Find.Execute("27+7=", $true, $true, $false, $false, $false, $true, 1, $false, "82-32=", 2) | Out-Null
$d.Content.Find.Execute("27-26=", $true, $true, $false, $false, $false, $true, 1, $false, "69-36=", 2) | Out-Null
$d.Content.Find.Execute("5-3=", $true, $true, $false, $false, $false, $true, 1, $false, "70-35=", 2) | Out-Null
$d.Content.Find.Execute("81-36=", $true, $true, $false, $false, $false, $true, 1, $false, "76-3=", 2) | Out-Null
$d.Content.Find.Execute("3+4=", $true, $true, $false, $false, $false, $true, 1, $false, "89-35=", 2) | Out-Null
$d.Content.Find.Execute("62+2=", $true, $true, $false, $false, $false, $true, 1, $false, "56-32=", 2) | Out-Null
$d.Content.Find.Execute("43+27=", $true, $true, $false, $false, $false, $true, 1, $false, "96-54=", 2) | Out-Null
$d.Content.Find.Execute("76-24=", $true, $true, $false, $false, $false, $true, 1, $false, "93-47=", 2) | Out-Null
$d.Content.Find.Execute("99-36=", $true, $true, $false, $false, $false, $true, 1, $false, "26+8=", 2) | Out-Null
$d.Content.Find.Execute("28+14=", $true, $true, $false, $false, $false, $true, 1, $false, "47-41=", 2) | Out-Null
$d.Content.Find.Execute("69-61=", $true, $true, $false, $false, $false, $true, 1, $false, "48+6=", 2) | Out-Null
$d.Content.Find.Execute("55-36=", $true, $true, $false, $false, $false, $true, 1, $false, "17+57=", 2) | Out-Null
$d.Content.Find.Execute("99-4=", $true, $true, $false, $false, $false, $true, 1, $false, "0+47=", 2) | Out-Null
$d.Content.Find.Execute("25-24=", $true, $true, $false, $false, $false, $true, 1, $false, "32-3=", 2) | Out-Null
$d.Content.Find.Execute("6+18=", $true, $true, $false, $false, $false, $true, 1, $false, "79-70=", 2) | Out-Null
$d.Content.Find.Execute("74-23=", $true, $true, $false, $false, $false, $true, 1, $false, "14+5=", 2) | Out-Null
$d.Content.Find.Execute("94-13=", $true, $true, $false, $false, $false, $true, 1, $false, "11+29=", 2) | Out-Null
$d.Content.Find.Execute("98-70=", $true, $true, $false, $false, $false, $true, 1, $false, "80-68=", 2) | Out-Null
$d.Content.Find.Execute("55-32=", $true, $true, $false, $false, $false, $true, 1, $false, "70-24=", 2) | Out-Null
$d.Content.Find.Execute("17+2=", $true, $true, $false, $false, $false, $true, 1, $false, "85-60=", 2) | Out-Null
$d.Content.Find.Execute("83+9=", $true, $true, $false, $false, $false, $true, 1, $false, "77-27=", 2) | Out-Null
$d.Content.Find.Execute("93-71=", $true, $true, $false, $false, $false, $true, 1, $false, "77+18=", 2) | Out-Null
$d.Content.Find.Execute("52+28=", $true, $true, $false, $false, $false, $true, 1, $false, "85+14=", 2) | Out-Null
$d.Content.Find.Execute("77-12=", $true, $true, $false, $false, $false, $true, 1, $false, "49+16=", 2) | Out-Null
$d.Content.Find.Execute("69-34=", $true, $true, $false, $false, $false, $true, 1, $false, "88-84=", 2) | Out-Null
$d.Content.Find.Execute("29+4=", $true, $true, $false, $false, $false, $true, 1, $false, "2+61=", 2) | Out-Null
$d.Content.Find.Execute("48+48=", $true, $true, $false, $false, $false, $true, 1, $false, "57+40=", 2) | Out-Null
$d.Content.Find.Execute("50+32=", $true, $true, $false, $false, $false, $true, 1, $false, "82-35=", 2) | Out-Null
$d.Content.Find.Execute("37+45=", $true, $true, $false, $false, $false, $true, 1, $false, "29-10=", 2) | Out-Null
$d.Content.Find.Execute("99-66=", $true, $true, $false, $false, $false, $true, 1, $false, "61-55=", 2) | Out-Null
$d.Content.Find.Execute("17+60=", $true, $true, $false, $false, $false, $true, 1, $false, "23+6=", 2) | Out-Null
$d.Content.Find.Execute("57+5=", $true, $true, $false, $false, $false, $true, 1, $false, "61-30=", 2) | Out-Null
$d.Content.Find.Execute("10+50=", $true, $true, $false, $false, $false, $true, 1, $false, "42+26=", 2) | Out-Null
$d.Content.Find.Execute("82-81=", $true, $true, $false, $false, $false, $true, 1, $false, "69-30=", 2) | Out-Null
$d.Content.Find.Execute("75+9=", $true, $true, $false, $false, $false, $true, 1, $false, "73+1=", 2) | Out-Null
$d.Content.Find.Execute("49+1=", $true, $true, $false, $false, $false, $true, 1, $false, "98-10=", 2) | Out-Null
$d.Content.Find.Execute("24+22=", $true, $true, $false, $false, $false, $true, 1, $false, "4-2=", 2) | Out-Null
$d.Content.Find.Execute("20-2=", $true, $true, $false, $false, $false, $true, 1, $false, "68+0=", 2) | Out-Null
$d.Content.Find.Execute("82+9=", $true, $true, $false, $false, $false, $true, 1, $false, "76+4=", 2) | Out-Null
$d.Content.Find.Execute("44+41=", $true, $true, $false, $false, $false, $true, 1, $false, "56-19=", 2) | Out-Null
$d.Content.Find.Execute("18+64=", $true, $true, $false, $false, $false, $true, 1, $false, "67-44=", 2) | Out-Null
$d.Content.Find.Execute("84-49=", $true, $true, $false, $false, $false, $true, 1, $false, "89-10=", 2) | Out-Null
$d.Content.Find.Execute("42+22=", $true, $true, $false, $false, $false, $true, 1, $false, "77-59=", 2) | Out-Null
$d.Content.Find.Execute("48+11=", $true, $true, $false, $false, $false, $true, 1, $false, "99-48=", 2) | Out-Null
$d.Content.Find.Execute("87-38=", $true, $true, $false, $false, $false, $true, 1, $false, "37-11=", 2) | Out-Null
$d.Content.Find.Execute("84-55=", $true, $true, $false, $false, $false, $true, 1, $false, "27+8=", 2) | Out-Null
$d.Content.Find.Execute("41+3=", $true, $true, $false, $false, $false, $true, 1, $false, "84-81=", 2) | Out-Null
$d.Content.Find.Execute("26+37=", $true, $true, $false, $false, $false, $true, 1, $false, "94-70=", 2) | Out-Null
$d.Content.Find.Execute("54-8=", $true, $true, $false, $false, $false, $true, 1, $false, "57+6=", 2) | Out-Null
$d.Content.Find.Execute("62-10=", $true, $true, $false, $false, $false, $true, 1, $false, "34-32=", 2) | Out-Null
$d.Content.Find.Execute("85-81=", $true, $true, $false, $false, $false, $true, 1, $false, "4+14=", 2) | Out-Null
$d.Content.Find.Execute("51+20=", $true, $true, $false, $false, $false, $true, 1, $false, "87-42=", 2) | Out-Null
$d.Content.Find.Execute("81+8=", $true, $true, $false, $false, $false, $true, 1, $false, "20-0=", 2) | Out-Null
$d.Content.Find.Execute("40+23=", $true, $true, $false, $false, $false, $true, 1, $false, "68+27=", 2) | Out-Null
$d.Content.Find.Execute("8+23=", $true, $true, $false, $false, $false, $true, 1, $false, "18+36=", 2) | Out-Null
$d.Content.Find.Execute("53-24=", $true, $true, $false, $false, $false, $true, 1, $false, "2+49=", 2) | Out-Null
$d.Content.Find.Execute("96-17=", $true, $true, $false, $false, $false, $true, 1, $false, "19+0=", 2) | Out-Null
$d.Content.Find.Execute("53-21=", $true, $true, $false, $false, $false, $true, 1, $false, "75-13=", 2) | Out-Null
$d.Content.Find.Execute("54+29=", $true, $true, $false, $false, $false, $true, 1, $false, "57+7=", 2) | Out-Null
$d.Content.Find.Execute("43+12=", $true, $true, $false, $false, $false, $true, 1, $false, "48-42=", 2) | Out-Null
$d.Content.Find.Execute("45-10=", $true, $true, $false, $false, $false, $true, 1, $false, "30+46=", 2) | Out-Null
$d.Content.Find.Execute("51-20=", $true, $true, $false, $false, $false, $true, 1, $false, "40+3=", 2) | Out-Null
$d.Content.Find.Execute("71-25=", $true, $true, $false, $false, $false, $true, 1, $false, "27+70=", 2) | Out-Null
$d.Content.Find.Execute("37+10=", $true, $true, $false, $false, $false, $true, 1, $false, "50+5=", 2) | Out-Null
$d.Content.Find.Execute("97-30=", $true, $true, $false, $false, $false, $true, 1, $false, "95-51=", 2) | Out-Null
$d.Content.Find.Execute("17+16=", $true, $true, $false, $false, $false, $true, 1, $false, "13+40=", 2) | Out-Null
$d.Content.Find.Execute("12+61=", $true, $true, $false, $false, $false, $true, 1, $false, "45-41=", 2) | Out-Null
$d.Content.Find.Execute("50+33=", $true, $true, $false, $false, $false, $true, 1, $false, "56-32=", 2) | Out-Null
$d.Content.Find.Execute("10+53=", $true, $true, $false, $false, $false, $true, 1, $false, "90-35=", 2) | Out-Null
$d.Content.Find.Execute("27-2=", $true, $true, $false, $false, $false, $true, 1, $false, "59-14=", 2) | Out-Null
$d.Content.Find.Execute("51-34=", $true, $true, $false, $false, $false, $true, 1, $false, "85-29=", 2) | Out-Null
$d.Content.Find.Execute("25+36=", $true, $true, $false, $false, $false, $true, 1, $false, "59-28=", 2) | Out-Null
$d.Content.Find.Execute("25-17=", $true, $true, $false, $false, $false, $true, 1, $false, "54+26=", 2) | Out-Null
$d.Content.Find.Execute("86-34=", $true, $true, $false, $false, $false, $true, 1, $false, "82-39=", 2) | Out-Null
$d.Content.Find.Execute("14+58=", $true, $true, $false, $false, $false, $true, 1, $false, "46-23=", 2) | Out-Null
$d.Content.Find.Execute("29+65=", $true, $true, $false, $false, $false, $true, 1, $false, "4+13=", 2) | Out-Null
$d.Content.Find.Execute("12+87=", $true, $true, $false, $false, $false, $true, 1, $false, "84-81=", 2) | Out-Null
$d.Content.Find.Execute("3+76=", $true, $true, $false, $false, $false, $true, 1, $false, "71+13=", 2) | Out-Null
$d.Content.Find.Execute("4+30=", $true, $true, $false, $false, $false, $true, 1, $false, "43+56=", 2) | Out-Null
$d.Content.Find.Execute("79-34=", $true, $true, $false, $false, $false, $true, 1, $false, "94-37=", 2) | Out-Null
$d.Content.Find.Execute("63+25=", $true, $true, $false, $false, $false, $true, 1, $false, "15+29=", 2) | Out-Null
$d.Content.Find.Execute("67-31=", $true, $true, $false, $false, $false, $true, 1, $false, "1+45=", 2) | Out-Null
$d.Content.Find.Execute("17+51=", $true, $true, $false, $false, $false, $true, 1, $false, "33-25=", 2) | Out-Null
$d.Content.Find.Execute("72-64=", $true, $true, $false, $false, $false, $true, 1, $false, "28-19=", 2) | Out-Null
$d.Content.Find.Execute("50+46=", $true, $true, $false, $false, $false, $true, 1, $false, "96-45=", 2) | Out-Null
$d.Content.Find.Execute("82-60=", $true, $true, $false, $false, $false, $true, 1, $false, "58-33=", 2) | Out-Null
$d.Content.Find.Execute("28-1=", $true, $true, $false, $false, $false, $true, 1, $false, "94-75=", 2) | Out-Null
$d.Content.Find.Execute("21+26=", $true, $true, $false, $false, $false, $true, 1, $false, "8+55=", 2) | Out-Null
$d.Content.Find.Execute("20-17=", $true, $true, $false, $false, $false, $true, 1, $false, "59-59=", 2) | Out-Null
$d.Content.Find.Execute("49+32=", $true, $true, $false, $false, $false, $true, 1, $false, "90-30=", 2) | Out-Null
$d.Content.Find.Execute("17+48=", $true, $true, $false, $false, $false, $true, 1, $false, "71-52=", 2) | Out-Null
$d.Content.Find.Execute("38+6=", $true, $true, $false, $false, $false, $true, 1, $false, "8+2=", 2) | Out-Null
